$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "-"
$ws.Range("D4").Value = "MCT-2A-Eletrônica analóg. e de potência"
$ws.Range("F4").Value = "-"
$ws.Range("D6").Value = "MCT-2A-Eletrônica analóg. e de potência"

$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"

$ws.Range("E18").Value = "-"
$ws.Range("E19").Value = "-"

$ws.Range("B20").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("D20").Value = "-"

$ws.Range("B21").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
